$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for rows 2, 3 and 4 in the columns
# that actually differ between the rows: A, B, D, E, F, G, H, Q, R.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$orig = @{}
foreach ($r in 2..4) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $orig[$r] = $rowVals
}

# Apply a cyclic rotation of the row data:
#   new row2 = old row4
#   new row3 = old row2
#   new row4 = old row3
$mapping = @{ 2 = 4; 3 = 2; 4 = 3 }

foreach ($destRow in 2..4) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $orig[$srcRow][$col]
    }
}
